# Update data analysis report:
#  - translate the Chinese report text to English
#  - collapse each paragraph's multiple runs into a single plain run
#  - move the "_GoBack" bookmark from the first paragraph to the end of
#    the (now last, still last) paragraph

$d = $word.ActiveDocument

# The "_GoBack" bookmark is hidden (name starts with "_"), so it won't show
# up in Bookmarks.Count/enumeration, but it is still reachable by name.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

function Set-ParaText($para, [string]$text) {
    # Replace the paragraph's whole content (but not its paragraph mark)
    # with a single run of plain text. Using Delete()+InsertAfter() (rather
    # than a Range.Text assignment) guarantees the old runs are actually
    # removed and replaced by one fresh run -- Range.Text is a no-op when
    # the replacement text happens to equal the existing text, which would
    # otherwise leave the original multi-run split (and its rPr) in place.
    $r = $para.Range
    [void]$r.MoveEnd(1, -1)
    $r.Delete()
    $r.InsertAfter($text)
}

Set-ParaText $d.Paragraphs.Item(1) "Conclusion: Positive correlation, good correlation"
Set-ParaText $d.Paragraphs.Item(2) "A total of 14 sets of data"
Set-ParaText $d.Paragraphs.Item(3) "Metric5&6"
Set-ParaText $d.Paragraphs.Item(4) "Pearson correlation coefficient 0.2732"
Set-ParaText $d.Paragraphs.Item(5) "Spearman correlation coefficient 0.0990"

# Re-create "_GoBack" collapsed at the end of the last paragraph's text
# (right before its paragraph mark). Bookmarks.Add on a collapsed range
# sitting exactly at the end of the document's content gets mis-resolved
# by this host, so insert a throwaway character to anchor a tiny 1-char
# range, bookmark that, then delete the character -- the bookmark
# collapses down to the correct position and stays put.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $last.Range
[void]$r.MoveEnd(1, -1)
$r.Collapse(0)
$r.InsertAfter("X")
$anchor = $d.Range($r.Start, $r.Start + 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range($r.Start, $r.Start + 1).Delete()

Write-Output "ok"
